$wb = $excel.ActiveWorkbook

$wsInvoice = $wb.Worksheets.Item("Historical Invoice Template")
$wsInvoice.Range("B2").Value = "Invoice1497179"
$wsInvoice.Range("C2").Value = "Invoice1764610"
$wsInvoice.Range("E2").Value = "Edwardct2l"
$wsInvoice.Range("F2").Value = "Aviva4x68"

$wsPO = $wb.Worksheets.Item("Historical PO Template")
$wsPO.Range("B2").Value = "Invoice1764610"
$wsPO.Range("C2").Value = "Invoice1497179"
$wsPO.Range("E2").Value = "Edwardct2l"
$wsPO.Range("F2").Value = "Aviva4x68"
